# Updates cryptos list prices (column D) and 1h volume change percentages
# (column E) for the rows whose underlying market data changed, matching
# the "Updated cryptos list ... with GitHub Actions" commit.
#
# The source cells are text (e.g. "57.762.00", "  -3.43%  ") rather than
# numbers -- some even use '.' as a thousands separator, which Excel would
# otherwise happily reinterpret as a number/date. To keep them as literal
# text (like the original file) we flip the cell to the "Text" number
# format ("@") just long enough to assign the value, then call
# ClearFormats() so the cell's style/format reverts to its original
# (default) state -- only the cell's text content ends up changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "57.762.00" },
    @{ Cell = "E2"; Value = "  -3.43%  " },
    @{ Cell = "D3"; Value = "2.273.68" },
    @{ Cell = "E3"; Value = "  -3.79%  " },
    @{ Cell = "E4"; Value = "  +0.06%  " },
    @{ Cell = "D5"; Value = "529.77" },
    @{ Cell = "E5"; Value = "  -5.36%  " },
    @{ Cell = "D6"; Value = "130.45" },
    @{ Cell = "E6"; Value = "  -2.18%  " },
    @{ Cell = "E7"; Value = "  +0.07%  " },
    @{ Cell = "D8"; Value = "0.582" },
    @{ Cell = "E8"; Value = "  -0.24%  " },
    @{ Cell = "D9"; Value = "2.273.17" },
    @{ Cell = "E9"; Value = "  -3.72%  " },
    @{ Cell = "D10"; Value = "0.0987" },
    @{ Cell = "E10"; Value = "  -6.01%  " },
    @{ Cell = "E11"; Value = "  -3.42%  " },
    @{ Cell = "E12"; Value = "  -0.31%  " },
    @{ Cell = "E13"; Value = "  -4.04%  " },
    @{ Cell = "D14"; Value = "23.39" },
    @{ Cell = "E14"; Value = "  -3.37%  " },
    @{ Cell = "D15"; Value = "2.679.44" },
    @{ Cell = "E15"; Value = "  -3.84%  " },
    @{ Cell = "D16"; Value = "57.749.18" },
    @{ Cell = "E16"; Value = "  -3.37%  " },
    @{ Cell = "E17"; Value = "  -4.63%  " },
    @{ Cell = "D18"; Value = "2.271.37" },
    @{ Cell = "E18"; Value = "  -4.38%  " },
    @{ Cell = "D19"; Value = "10.47" },
    @{ Cell = "E19"; Value = "  -5.06%  " },
    @{ Cell = "D20"; Value = "4.19" },
    @{ Cell = "E20"; Value = "  -6.05%  " },
    @{ Cell = "D21"; Value = "310.76" },
    @{ Cell = "E21"; Value = "  -2.88%  " },
    @{ Cell = "E22"; Value = "  -4.45%  " },
    @{ Cell = "E23"; Value = "  -0.08%  " },
    @{ Cell = "D24"; Value = "62.42" },
    @{ Cell = "E24"; Value = "  -2.63%  " },
    @{ Cell = "D25"; Value = "0.167" },
    @{ Cell = "E25"; Value = "  -3.39%  " },
    @{ Cell = "E26"; Value = "  +0.05%  " },
    @{ Cell = "D27"; Value = "7.94" },
    @{ Cell = "E27"; Value = "  -5.28%  " },
    @{ Cell = "E28"; Value = "  -7.29%  " },
    @{ Cell = "D29"; Value = "170.48" },
    @{ Cell = "E29"; Value = "  -0.02%  " },
    @{ Cell = "E30"; Value = "  -6.59%  " },
    @{ Cell = "E31"; Value = "  -5.71%  " },
    @{ Cell = "D32"; Value = "5.71" },
    @{ Cell = "E32"; Value = "  -5.88%  " },
    @{ Cell = "D33"; Value = "1.03" },
    @{ Cell = "E33"; Value = "  -6.17%  " },
    @{ Cell = "E36"; Value = "  -2.21%  " },
    @{ Cell = "E37"; Value = "  +0.03%  " },
    @{ Cell = "E38"; Value = "  -7.00%  " },
    @{ Cell = "E39"; Value = "  -5.76%  " },
    @{ Cell = "D40"; Value = "38.22" },
    @{ Cell = "E40"; Value = "  -1.02%  " },
    @{ Cell = "E41"; Value = "  -6.67%  " },
    @{ Cell = "D42"; Value = "140.72" },
    @{ Cell = "E42"; Value = "  -2.85%  " },
    @{ Cell = "D43"; Value = "284.62" },
    @{ Cell = "E43"; Value = "  -10.53%  " },
    @{ Cell = "E44"; Value = "  -3.43%  " },
    @{ Cell = "E45"; Value = "  -2.10%  " },
    @{ Cell = "D46"; Value = "0.0492" },
    @{ Cell = "E46"; Value = "  -3.31%  " },
    @{ Cell = "D47"; Value = "0.548" },
    @{ Cell = "E47"; Value = "  -3.06%  " },
    @{ Cell = "D48"; Value = "17.96" },
    @{ Cell = "E48"; Value = "  -7.41%  " },
    @{ Cell = "E49"; Value = "  -3.90%  " },
    @{ Cell = "D50"; Value = "10.93" },
    @{ Cell = "E50"; Value = "  -1.21%  " },
    @{ Cell = "D51"; Value = "4.64" },
    @{ Cell = "E51"; Value = "  -0.41%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
